$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITERACION2")

# --- Update requirement text cells (2nd iteration revision) ---
# Order matches the original authoring sequence so new shared-string
# entries are appended in the same order as the target workbook.

# RQ09 (row 12): expand the description of the "Realizar denuncia" requirement
$ws.Range("E12").Value = "El sistema debe permitir ingresar los datos (archivos como imágenes, video, audio, descripcion) de una denuncia y registrarlos. "

# RQ01 (row 4): "codigo de placa" -> "numero de placa"
$ws.Range("D4").Value = "Ingresar numero de placa"
$ws.Range("E4").Value = "El sistema debe permitir ingresar en un campo de texto el numero de placa"

# RQ02 (row 5): "codigo de placa" -> "numero de placa"
$ws.Range("D5").Value = "Buscar por numero de placa"

# RQ05 (row 8): "codigo de placa" -> "numero de placa"
$ws.Range("D8").Value = "Detectar numero de placa en imagen"

$ws.Range("E5").Value = "El sistema debe permitir buscar por numero de placa"

# RQ15 (row 17): "codigo de placa" -> "numero de placa"
$ws.Range("D17").Value = "Buscar por numero de placa"
$ws.Range("E17").Value = "El sistema debe permitir buscar por numero de placa"

# RQ19 (row 21): "codigo de placa" -> "numero de placa"
$ws.Range("D21").Value = "Buscar por numero de placa"
$ws.Range("E21").Value = "El sistema debe permitir buscar por numero de placa"

# --- Update the view state: scroll position + active selection ---
$ws.Range("E16").Select()
